$d = $word.ActiveDocument
$apos = [char]0x2019

# ---------------------------------------------------------------------------
# 1. Remove the _GoBack bookmark that currently sits right after "Work Log".
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2. Replace the placeholder "sgrdg" text with the full sentence that
#    describes the work completed on April 30th, 2017. We build up the
#    paragraph's text in a single Find/Replace so every character keeps the
#    Times New Roman formatting that was already present on that run.
# ---------------------------------------------------------------------------
$marker = "Annyang" + $apos + "s"
$fullSentence = $marker + " voice commands now add and removed items from the to-do list.  Wesley and Michael have been working on pulls and pushes to mongoDB through the client.  We are successfully able to read and write user information to the database and have it saved.  We started on the required documents to include with the project (README, Technologies, Guide)."

$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("sgrdg", $true, $false, $false, $false, $false, $true, 1, $false, $fullSentence, 2)

# ---------------------------------------------------------------------------
# 3. Split the merged run into the pieces that existed in the target
#    revision (a separate run for "Annyang's" and a separate run for
#    "mongoDB") by briefly toggling a character property, which forces the
#    engine to break the run apart while copying the existing rPr (so the
#    Times New Roman formatting survives) onto each new piece.
# ---------------------------------------------------------------------------
function Split-Run($range) {
    $range.Bold = 1
    $range.Bold = 0
}

$findMarker = $d.Content.Find
$findMarker.ClearFormatting()
$findMarker.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Split-Run $findMarker.Parent

$findMongo = $d.Content.Find
$findMongo.ClearFormatting()
$findMongo.Execute("mongoDB", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Split-Run $findMongo.Parent

# ---------------------------------------------------------------------------
# 4. Re-insert the _GoBack bookmark, collapsed at the very end of the
#    paragraph we just edited. Bookmarks.Add() on a truly-collapsed range
#    mis-anchors the bookmarkStart marker, so we add it over a one-character
#    placeholder and then delete the placeholder; the bookmark shrinks back
#    to a zero-width bookmark at the correct position.
# ---------------------------------------------------------------------------
$findEnd = $d.Content.Find
$findEnd.ClearFormatting()
$findEnd.Execute("Guide).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $findEnd.Parent.End

$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($endPos, $endPos + 1).Text = ""
